$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "second_map"
$ws.Range("C2").Value = 512
$ws.Range("D2").Value = 64
$ws.Range("E2").Value = "SNORLAX.png"

$ws.Range("E2").Select()
